$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.939.74"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "1.650.92"
$ws.Range("E3").Value = "  +2.86%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'214.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("E6").Value = "  +2.62%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +2.68%  "

$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").Value = "'20.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.75%  "

$ws.Range("E11").Value = "  +2.49%  "

$ws.Range("D12").Value = "1.884.91"
$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("D13").Value = "1.645.36"
$ws.Range("E13").Value = "  +2.47%  "

$ws.Range("E14").Value = "  +1.83%  "

$ws.Range("E15").Value = "  +2.65%  "

$ws.Range("E16").Value = "  +2.76%  "

$ws.Range("D17").Value = "26.947.18"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "'235.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "'7.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").Value = "  +3.28%  "

$ws.Range("E23").Value = "  +3.94%  "

$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").Value = "'145.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").Value = "'7.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "

$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'15.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.63%  "

$ws.Range("D30").Value = "'0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").Value = "1.557.14"
$ws.Range("E32").Value = "  +4.40%  "

$ws.Range("D33").Value = "'3.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.62%  "

$ws.Range("E34").Value = "  +4.79%  "

$ws.Range("D35").Value = "'1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.54%  "

$ws.Range("D36").Value = "'2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").Value = "'0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.11%  "

$ws.Range("D38").Value = "'0.893"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.87%  "

$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "'65.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.41%  "

$ws.Range("E43").Value = "  +2.45%  "

$ws.Range("D44").Value = "1.791.74"
$ws.Range("E44").Value = "  +2.79%  "

$ws.Range("D45").Value = "'0.776"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("D46").Value = "'0.925"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.86%  "

$ws.Range("D47").Value = "'90.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").Value = "'0.0990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.07%  "

$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("E51").Value = "  +2.33%  "
